$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 25, shifting existing rows 25-30 down to 26-31
$ws.Rows.Item(25).Insert()

# Populate the new row 25 with data (values carried from the cell above for
# unchanged columns, per the diff)
$ws.Range("A25").Value = 10
$ws.Range("B25").Value = "Vega Modelo de Temuco"
$ws.Range("C25").Value = "La Araucanía"
$ws.Range("D25").Value = 44474
$ws.Range("E25").Value = 9
$ws.Range("F25").Value = 100112026
$ws.Range("G25").Value = "Haba"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 30
$ws.Range("K25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = 10000
$ws.Range("N25").Value = "$/saco 25 kilos"
$ws.Range("O25").Value = "Provincia de Limarí"
$ws.Range("P25").Value = 400
$ws.Range("Q25").Value = 25
$ws.Range("R25").Value = "Hortaliza"
